$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row (row 30) - "DRAIAM071" / customer-care web form script.
$tcid        = "DRAIAM071"
$jira        = "OPQA-5168 || OPQA-5227"
$description = "Verify that the web form provided to user should be application specific with following required fields`n1.Name 2.Organization 3.Contact details (email, telephone) 4.Issue Category 5.Country`n6.Description of issue ( a free form text box where a user can describe why they are contacting support) ||`nVerify that the web form provided to user will be application specific with following required fields (Name, Organization, email, telephone, Issue Category, Country, Description of issue)"
$runmode     = "Y"

# Mirror the formatting of the row directly above (thin borders all round,
# no wrap) and reuse the wrapped-text look from column C of row 28 for the
# long description cell - matches the rest of the sheet's existing styles.
$ws.Range("A29:E29").Copy($ws.Range("A30:E30"))
$ws.Range("C28").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("A30").Value = $tcid
$ws.Range("B30").Value = $jira
$ws.Range("C30").Value = $description
$ws.Range("D30").Value = $runmode
$ws.Range("E30").Value = ""

$ws.Rows.Item(30).RowHeight = 75

[void]$ws.Range("C30").Select()
